# Update column F (dSF) values for specific rows, per "repull data, push all
# data, mean calculation" commit. Column F was previously a straight copy of
# column E (dS0); this re-pull adjusts several rows' dSF figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -2
    5  = -4
    8  = -2
    10 = -3
    11 = 1
    13 = -1
    18 = -3
    19 = -4
    25 = -7
    30 = 1
    31 = 4
    32 = 5
    33 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
